$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the two formats used on the sheet (plain wrapped = row 1 style,
#     highlighted wrapped = row 2 A:F style) onto scratch rows far below the
#     data, then wipe the whole used range so stale formatting from rows that
#     change role (e.g. old test-name rows becoming plain step rows) does not
#     bleed through.
$ws.Range("A1:F1").Copy()
$ws.Range("A200:F200").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A201:F201").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1:F17").Clear()

# --- Re-apply the plain style to the header row and to every step-row
#     Step/Description/Expected triple (cols C:E), and the highlighted style
#     to every test-name row (cols A:F). Copy source widths are kept equal to
#     the paste target widths so PasteSpecial does not spill into extra cols.
$ws.Range("A200:F200").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C200:E200").Copy()
$ws.Range("C3:E3").PasteSpecial(-4122)
$ws.Range("C4:E4").PasteSpecial(-4122)
$ws.Range("C6:E6").PasteSpecial(-4122)
$ws.Range("C7:E7").PasteSpecial(-4122)
$ws.Range("C9:E9").PasteSpecial(-4122)
$ws.Range("C10:E10").PasteSpecial(-4122)
$ws.Range("C11:E11").PasteSpecial(-4122)
$ws.Range("C12:E12").PasteSpecial(-4122)
$ws.Range("C13:E13").PasteSpecial(-4122)
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Range("C16:E16").PasteSpecial(-4122)
$ws.Range("C17:E17").PasteSpecial(-4122)
$ws.Range("C18:E18").PasteSpecial(-4122)
$ws.Range("C20:E20").PasteSpecial(-4122)
$ws.Range("C21:E21").PasteSpecial(-4122)
$ws.Range("C22:E22").PasteSpecial(-4122)
$ws.Range("C24:E24").PasteSpecial(-4122)
$ws.Range("C25:E25").PasteSpecial(-4122)
$ws.Range("C26:E26").PasteSpecial(-4122)
$ws.Range("C27:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A201:F201").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A8:F8").PasteSpecial(-4122)
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A23:F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Remove the scratch rows so they do not show up in the saved sheet.
$ws.Range("A200:F201").Clear()

# --- Write the cell values for the full new layout (rows 1-27).
$ws.Cells.Item(1,1).Value = "Test Name"
$ws.Cells.Item(1,2).Value = "Test Description"
$ws.Cells.Item(1,3).Value = "Step Name"
$ws.Cells.Item(1,4).Value = "Step Description"
$ws.Cells.Item(1,5).Value = "Expected result"
$ws.Cells.Item(1,6).Value = "Actual result"

$ws.Cells.Item(2,1).Value = "01_SignIn_A_POS"
$ws.Cells.Item(2,2).Value = "Sign in with valid user information into the website"
$ws.Cells.Item(2,3).Value = "Precondition"
$ws.Cells.Item(2,4).Value = "Website is accessible.`nUser with login information exists"
$ws.Cells.Item(2,5).Value = "Environment and data available"
$ws.Cells.Item(2,6).Value = ""

$ws.Cells.Item(3,3).Value = "Step 1"
$ws.Cells.Item(3,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(3,5).Value = "Browser starts on specified page"

$ws.Cells.Item(4,3).Value = "Step 2"
$ws.Cells.Item(4,4).Value = "Use valid login information to sign into the website"
$ws.Cells.Item(4,5).Value = "User logged in"

$ws.Cells.Item(5,1).Value = "02_SignIn_A_NEG"
$ws.Cells.Item(5,2).Value = "Sign in with invalid user information into the website"
$ws.Cells.Item(5,3).Value = "Precondition"
$ws.Cells.Item(5,4).Value = "Website is accessible.`nUser with login information exists"
$ws.Cells.Item(5,5).Value = "Environment and data available"
$ws.Cells.Item(5,6).Value = ""

$ws.Cells.Item(6,3).Value = "Step 1"
$ws.Cells.Item(6,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(6,5).Value = "Browser starts on specified page"

$ws.Cells.Item(7,3).Value = "Step 2"
$ws.Cells.Item(7,4).Value = "Use invalid login information to sign into the website"
$ws.Cells.Item(7,5).Value = "User login refused and redirected to `"SIGN-ON`" page"

$ws.Cells.Item(8,1).Value = "03_FlightReservation_A_POS"
$ws.Cells.Item(8,2).Value = ""
$ws.Cells.Item(8,3).Value = "Precondition"
$ws.Cells.Item(8,4).Value = "Website is accessible.`nUser with login information exists"
$ws.Cells.Item(8,5).Value = "Environment and data available"
$ws.Cells.Item(8,6).Value = ""

$ws.Cells.Item(9,3).Value = "Step 1"
$ws.Cells.Item(9,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(9,5).Value = "Browser starts on specified page"

$ws.Cells.Item(10,3).Value = "Step 2"
$ws.Cells.Item(10,4).Value = "Use valid login information to sign into the website"
$ws.Cells.Item(10,5).Value = "User logged in"

$ws.Cells.Item(11,3).Value = "Step 3"
$ws.Cells.Item(11,4).Value = "Fill out `"Flight Finder`" form and continue to next page with `"CONTINUE`" button under the form"
$ws.Cells.Item(11,5).Value = "Selected options applied into next form"

$ws.Cells.Item(12,3).Value = "Step 4"
$ws.Cells.Item(12,4).Value = "Select any DEPART and RETURN options in `"Select Flight`" form and continue to next page with `"CONTINUE`" button under the form"
$ws.Cells.Item(12,5).Value = "Selected options applied into next form"

$ws.Cells.Item(13,3).Value = "Step 5"
$ws.Cells.Item(13,4).Value = "Fill out `"Book a Flight`" form and continue to next page with `"SECURE PURCHASE`" button under the form"
$ws.Cells.Item(13,5).Value = "Selected options applied into receipt"

$ws.Cells.Item(14,1).Value = "01_RegisterNewUser_M_POS"
$ws.Cells.Item(14,2).Value = ""
$ws.Cells.Item(14,3).Value = "Precondition"
$ws.Cells.Item(14,4).Value = "Website is accessible."
$ws.Cells.Item(14,5).Value = "Environment and data available"
$ws.Cells.Item(14,6).Value = ""

$ws.Cells.Item(15,3).Value = "Step 1"
$ws.Cells.Item(15,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(15,5).Value = "Browser starts on specified page"

$ws.Cells.Item(16,3).Value = "Step 2"
$ws.Cells.Item(16,4).Value = "Open registration form by clicking `"REGISTER`" button in top menu"
$ws.Cells.Item(16,5).Value = "Registration form opens"

$ws.Cells.Item(17,3).Value = "Step 3"
$ws.Cells.Item(17,4).Value = "Fill in registration form and submit it with `"SUBMIT`" button bellow the form"
$ws.Cells.Item(17,5).Value = "Registration successful, new user is created"

$ws.Cells.Item(18,3).Value = "Step 4"
$ws.Cells.Item(18,4).Value = "Sign in using new user information"
$ws.Cells.Item(18,5).Value = "User logged in"

$ws.Cells.Item(19,1).Value = "02_SignOut_M_POS"
$ws.Cells.Item(19,2).Value = ""
$ws.Cells.Item(19,3).Value = "Precondition"
$ws.Cells.Item(19,4).Value = "Website is accessible.`nUser is logged into the website"
$ws.Cells.Item(19,5).Value = "Environment and data available"
$ws.Cells.Item(19,6).Value = ""

$ws.Cells.Item(20,3).Value = "Step 1"
$ws.Cells.Item(20,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(20,5).Value = "Browser starts on specified page"

$ws.Cells.Item(21,3).Value = "Step 2"
$ws.Cells.Item(21,4).Value = "Use valid login information to sign into the website"
$ws.Cells.Item(21,5).Value = "User logged in"

$ws.Cells.Item(22,3).Value = "Step 3"
$ws.Cells.Item(22,4).Value = "Sign out using `"SIGN-OFF`" button in top menu"
$ws.Cells.Item(22,5).Value = "User logged out"

$ws.Cells.Item(23,1).Value = "03_CancelAllReservations_M_POS"
$ws.Cells.Item(23,2).Value = ""
$ws.Cells.Item(23,3).Value = "Precondition"
$ws.Cells.Item(23,4).Value = "Website is accessible.`nUser with login information exists.`nUser has existing reservation."
$ws.Cells.Item(23,5).Value = "Environment and data available"
$ws.Cells.Item(23,6).Value = ""

$ws.Cells.Item(24,3).Value = "Step 1"
$ws.Cells.Item(24,4).Value = "Start Chrome browser on `"http://newtours.demoaut.com`" website"
$ws.Cells.Item(24,5).Value = "Browser starts on specified page"

$ws.Cells.Item(25,3).Value = "Step 2"
$ws.Cells.Item(25,4).Value = "Use valid login information to sign into the website"
$ws.Cells.Item(25,5).Value = "User logged in"

$ws.Cells.Item(26,3).Value = "Step 3"
$ws.Cells.Item(26,4).Value = "Open Itinerary page with `"ITINERARY`" button in top menu"
$ws.Cells.Item(26,5).Value = "Itinerary page opens with existing reservation"

$ws.Cells.Item(27,3).Value = "Step 4"
$ws.Cells.Item(27,4).Value = "Cancel all reservations with `"CANCEL ALL RESERVATIONS`" button at the bottom of the page"
$ws.Cells.Item(27,5).Value = "Reservations are canceled and not active"

# --- Row heights
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 90
$ws.Rows.Item(13).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 45
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 75
$ws.Rows.Item(24).RowHeight = 45
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 45
$ws.Rows.Item(27).RowHeight = 60

# --- Selection matches the saved cursor position in the target file
$ws.Range("D23").Select()

